# Automatic update of files.
#
# 1) Column C ("Förändrad") is refreshed from 45184 to 45186 for every
#    data row (2 through 347).
# 2) Every HYPERLINK() formula in columns S,T,U,V,W,X,Y (rows 2-20) gets a
#    second "friendly name" argument equal to the row's "Beteckning"
#    (column A) value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow  = 347
$hyperlinkCols = @("S", "T", "U", "V", "W", "X", "Y")

for ($row = $firstRow; $row -le $lastRow; $row++) {

    # --- 1) Bump the "Förändrad" date serial in column C (45184 -> 45186) ---
    $ws.Range("C$row").Value = 45186

    # --- 2) Add the friendly-name argument to HYPERLINK formulas ---
    $name = $ws.Range("A$row").Value()

    foreach ($col in $hyperlinkCols) {
        $cell = $ws.Range("$col$row")
        $oldFormula = $cell.Formula()

        if ($oldFormula -and $oldFormula.StartsWith("=HYPERLINK(")) {
            # Only touch formulas that don't already carry a second argument.
            $openParen = $oldFormula.IndexOf("(")
            $closeParen = $oldFormula.LastIndexOf(")")
            $argsText = $oldFormula.Substring($openParen + 1, $closeParen - $openParen - 1)

            if ($argsText.IndexOf(",") -lt 0) {
                $q1 = $argsText.IndexOf('"')
                $q2 = $argsText.IndexOf('"', $q1 + 1)
                $url = $argsText.Substring($q1 + 1, $q2 - $q1 - 1)

                $newFormula = '=HYPERLINK("' + $url + '", "' + $name + '")'
                $cell.Formula = $newFormula
            }
        }
    }
}
